# before.xlsx -> after.xlsx
# Changes per the commit diff:
#  1. Strip the leading "> " marker from the 5 pathway-name header cells
#     (A1:E1). Note C1 ("EGF Signaling Pathway") only loses the ">" and
#     keeps a leading space, matching the diff exactly.
#  2. Move the sheet's active cell/selection from E2 to E1.
#  3. Reduce the window's tab-ratio (tab bar vs. horizontal scrollbar
#     split) from 991/1000 to 500/1000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Signaling of Hepatocyte Growth Factor Receptor"
$ws.Range("B1").Value = "TGF-beta signaling pathway"
$ws.Range("C1").Value = " EGF Signaling Pathway"
$ws.Range("D1").Value = "EPO Signaling Pathway"
$ws.Range("E1").Value = "GAS6 Signaling Pathway"

# Move the selection/active cell from E2 to E1.
[void]$ws.Range("E1").Select()

# Shrink the tab ratio (window chrome state; matches the authored diff's
# bookViews/workbookView@tabRatio 991 -> 500).
$excel.ActiveWindow.TabRatio = 0.5
